# fixed bug in delay_ms function: desired Tloop was entered in ms (1000)
# instead of seconds (0.00434 s == 4.34 ms @ 115200 baud, 0.5 bit-times).
# Also add the reference calc for the delay_ms timing constant in B12.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Sheet1!C9 "desired Tloop (ms)" - correct the units bug
$ws1.Range("C9").Value = 0.00434

# Sheet1!B12 - new helper calc: 1000 * 0.5 bit / 115200 baud
$ws1.Range("B12").Formula = "=1000*0.5/115200"
